$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Add the four new systematic-uncertainty component headers (M1:P1) first so the
# shared-string table gains "syst0_c".."syst3_c" before "syst_tot" is appended.
$ws.Range("M1").Value = "syst0_c"
$ws.Range("N1").Value = "syst1_c"
$ws.Range("O1").Value = "syst2_c"
$ws.Range("P1").Value = "syst3_c"

# Style the new header cells: Calibri 11, centered horizontally, top-aligned vertically.
$newHeaders = $ws.Range("M1:P1")
$newHeaders.Font.Name = "Calibri"
$newHeaders.Font.Size = 11
$newHeaders.HorizontalAlignment = -4108
$newHeaders.VerticalAlignment = -4160

# Rename existing "syst_u" header (column H) to "syst_tot".
$ws.Range("H1").Value = "syst_tot"

# --- Data rows (2:7) : new syst0_c / syst1_c / syst2_c / syst3_c columns ---
$ws.Range("M2").Value = 0.0002
$ws.Range("N2").Value = 0.006999999999999999
$ws.Range("O2").Value = 0.0011
$ws.Range("P2").Value = 0.0016

$ws.Range("M3").Value = 0.0003
$ws.Range("N3").Value = 0.006999999999999999
$ws.Range("O3").Value = 0.0009
$ws.Range("P3").Value = 0.0019

$ws.Range("M4").Value = 0.0003
$ws.Range("N4").Value = 0.006999999999999999
$ws.Range("O4").Value = 0.0019
$ws.Range("P4").Value = 0.0026

$ws.Range("M5").Value = 0.0008
$ws.Range("N5").Value = 0.006999999999999999
$ws.Range("O5").Value = 0.004699999999999999
$ws.Range("P5").Value = 0.0033

$ws.Range("M6").Value = 0.0009
$ws.Range("N6").Value = 0.006999999999999999
$ws.Range("O6").Value = 0.004
$ws.Range("P6").Value = 0.0025

$ws.Range("M7").Value = 0.001
$ws.Range("N7").Value = 0.006999999999999999
$ws.Range("O7").Value = 0.004500000000000001
$ws.Range("P7").Value = 0.0025

# --- Selection / view state, matching the saved workbook's cursor position ---
$ws.Range("H12").Select()
